$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 598; this pushes the former rows 598-669
# down to 599-670 (dimension grows from A1:T669 to A1:T670).
$ws.Rows.Item(598).Insert()

# Populate the newly inserted row 598 with the new weekly record.
$ws.Range("A598").Value2 = 9
$ws.Range("B598").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C598").Value2 = "Metropolitana"
$ws.Range("D598").Value2 = 45077
$ws.Range("E598").Value2 = 13
$ws.Range("F598").Value2 = "Fruta"
$ws.Range("G598").Value2 = 100108
$ws.Range("H598").Value2 = "Tropicales y subtropicales"
$ws.Range("I598").Value2 = 100108002
$ws.Range("J598").Value2 = "Mango"
$ws.Range("K598").Value2 = "Sin especificar"
$ws.Range("L598").Value2 = "Primera"
$ws.Range("M598").Value2 = 580
$ws.Range("N598").Value2 = 7500
$ws.Range("O598").Value2 = 8500
$ws.Range("P598").Value2 = 7983
$ws.Range("Q598").Value2 = "$/bandeja 4 kilos"
$ws.Range("R598").Value2 = "Perú"
$ws.Range("S598").Value2 = 1996
$ws.Range("T598").Value2 = 4
